# The table of lat/long/address/image rows accumulated a long run of
# near-duplicate rows (Picture 4 .. Picture 20) all pointing at the same
# "Meadway Court" address/photo while the geocoder/map-matching settled
# on a final position. Collapse that run down to just the last
# (converged) sample, which is the Picture 21 row.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 is the header, row 2 = Picture 1, ... row 5 = Picture 4, ...
# row 21 = Picture 20, row 22 = Picture 21 (last / kept row).
# Deleting row index 5 repeatedly removes rows 5..21 (Pictures 4-20),
# leaving the final row (Picture 21) intact as the new last row.
$firstToDrop = 5
$lastToDrop = 21
$dropCount = $lastToDrop - $firstToDrop + 1

for ($i = 0; $i -lt $dropCount; $i++) {
    $t.Rows.Item($firstToDrop).Delete()
}
